$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- 1. "Bottom 5 locations:" -> "Bottom 5 states:" --------------------
$tb = $s.Shapes.Item(2)
$tr = $tb.TextFrame.TextRange
$sub = $tr.Characters(83, 19)
$sub.Text = "Bottom 5 states:"

# --- 2. Reposition / resize the state map picture -----------------------
$pic = $s.Shapes.Item(3)
$pic.Left = 332.9617614746094
$pic.Top = 51.882362365722656
$pic.Width = 618.1033325195312
$pic.Height = 509.4559326171875

# --- 3. Add a right-brace shape ------------------------------------------
$brace = $s.Shapes.AddShape(32, 100, 100, 50, 50)
$brace.Name = "Right Brace 2"
$brace.Line.Weight = 2.25
$brace.TextFrame.VerticalAnchor = 3
$brace.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$brace.Left = 209.6470947265625
$brace.Top = 160.94119262695312
$brace.Width = 26.470552444458008
$brace.Height = 109.05882263183594

# --- 4. Add the "Top 5 make up 32 % of all breweries" textbox -----------
$callout = $s.Shapes.AddTextbox(1, 100, 100, 200, 100)
$callout.TextFrame.WordWrap = -1
$callout.TextFrame.AutoSize = 1
$callout.Fill.Visible = 0
$callout.TextFrame.TextRange.Text = "Top 5 make up 32 % of all breweries"
$callout.TextFrame.TextRange.Font.Bold = -1
$callout.Left = 243.18984985351562
$callout.Top = 179.11898803710938
$callout.Width = 113.63370513916016
$callout.Height = 72.70315551757812
